# Apply updated cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "98.727.68"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3
$ws.Range("D3").Value = "3.342.11"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.08"
$ws.Range("E5").Value = "  -2.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "643.97"
$ws.Range("E6").Value = "  +1.23%  "

# Row 7
$ws.Range("E7").Value = "  +13.24%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.431"
$ws.Range("E8").Value = "  +9.42%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.10"
$ws.Range("E9").Value = "  +26.32%  "

# Row 10
$ws.Range("E10").Value = "  +0.01%  "

# Row 11
$ws.Range("D11").Value = "3.339.55"
$ws.Range("E11").Value = "  -0.91%  "

# Row 12
$ws.Range("E12").Value = "  +3.21%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.72"
$ws.Range("E13").Value = "  +20.79%  "

# Row 14
$ws.Range("E14").Value = "  +8.29%  "

# Row 15
$ws.Range("D15").Value = "98.438.83"
$ws.Range("E15").Value = "  -0.18%  "

# Row 16
$ws.Range("D16").Value = "3.967.74"
$ws.Range("E16").Value = "  -0.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.53"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18
$ws.Range("D18").Value = "3.342.09"
$ws.Range("E18").Value = "  -0.24%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.24"
$ws.Range("E19").Value = "  +17.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.74"
$ws.Range("E20").Value = "  +10.80%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "539.42"
$ws.Range("E21").Value = "  +9.23%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.53"
$ws.Range("E22").Value = "  -1.92%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.15"
$ws.Range("E23").Value = "  +8.73%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.446"
$ws.Range("E24").Value = "  +57.14%  "

# Row 25
$ws.Range("E25").Value = "  -4.90%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "100.94"
$ws.Range("E26").Value = "  +13.57%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.17"
$ws.Range("E27").Value = "  +7.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.54"
$ws.Range("E28").Value = "  +3.92%  "

# Row 29
$ws.Range("D29").Value = "3.513.25"

# Row 30
$ws.Range("E30").Value = "  +16.92%  "

# Row 31
$ws.Range("E31").Value = "  +0.00%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.05"
$ws.Range("E32").Value = "  +15.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.190"
$ws.Range("E33").Value = "  -3.54%  "

# Row 34
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.20"
$ws.Range("E35").Value = "  +5.02%  "

# Row 36
$ws.Range("E36").Value = "  +12.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.62"
$ws.Range("E37").Value = "  +2.72%  "

# Row 38
$ws.Range("E38").Value = "  +3.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.155"
$ws.Range("E39").Value = "  +2.95%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "522.45"
$ws.Range("E40").Value = "  +2.75%  "

# Row 41
$ws.Range("E41").Value = "  -0.54%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.89"
$ws.Range("E42").Value = "  -0.72%  "

# Row 43
$ws.Range("E43").Value = "  +2.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.813"
$ws.Range("E44").Value = "  +5.17%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0401"
$ws.Range("E45").Value = "  +24.45%  "

# Row 46
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.21"
$ws.Range("E46").Value = "  -4.41%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.04%  "

# Row 48
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.84"
$ws.Range("E48").Value = "  +20.04%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.04"
$ws.Range("E49").Value = "  +4.58%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "163.61"
$ws.Range("E50").Value = "  +1.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.94"
$ws.Range("E51").Value = "  +7.35%  "
